$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the three "SFONN01" tag values to "MFASS01" and flag their INFILE column as active (1)
$ws.Range("E15").Value = "B017237ABO19_1_J_SIN_PMSL005_MFASS01"
$ws.Range("G15").Value = 1

$ws.Range("E14").Value = "B012935ABO13_1_J_SIN_PMSL001_MFASS01"
$ws.Range("G14").Value = 1

$ws.Range("E13").Value = "B012719ABO13_1_J_SIN_PMSL008_MFASS01"
$ws.Range("G13").Value = 1

# Update the view: scroll position and active/selected cell
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E14").Select()

$wb.Save()
